# edit.ps1 - applies the resume edits described by the target diff:
#   1. Skills bullet: "Visual Basic" -> "HTML, CSS, Javascript"
#   2. New bullet "Performed integrated testing of the entire system."
#      inserted after the "...I/O systems." bullet.
#   3. The "Used CAD software ... Arduino Uno microcontroller." bullet is
#      split into two bullets:
#        "Applied AutoCAD to create 3D printing model prototypes."
#        "Developed a Python-based GUI for the microcontroller."
#   4. New bullet "Navigated business landscapes with decisive action,
#      turning any obstacle into opportunities." inserted after the
#      "Built a loyal customer base ... for business thrived." bullet.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Programming & Data skills line: replace "Visual Basic" with the
#    new "HTML, CSS, Javascript" entry.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Visual Basic", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "HTML, CSS, Javascript", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Insert a new bullet after the "...controllers, and I/O systems."
#    line (same numbered-list paragraph formatting carries over
#    automatically via InsertParagraphAfter).
# ---------------------------------------------------------------------
$systemsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "and I/O systems\.") {
        $systemsPara = $p
        break
    }
}
$systemsPara.Range.InsertParagraphAfter()
$newPara = $systemsPara.Next()
$newPara.Range.InsertBefore("Performed integrated testing of the entire system.")

# ---------------------------------------------------------------------
# 3) Split the "Used CAD software..." bullet into two bullets.
# ---------------------------------------------------------------------
$cadPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Used CAD software") {
        $cadPara = $p
        break
    }
}
$cadRange = $cadPara.Range
$cadTextRange = $d.Range($cadRange.Start, $cadRange.End - 1)
$cadTextRange.Text = "Applied AutoCAD to create 3D printing model prototypes."

$cadPara.Range.InsertParagraphAfter()
$guiPara = $cadPara.Next()
$guiPara.Range.InsertBefore("Developed a Python-based GUI for the microcontroller.")

# ---------------------------------------------------------------------
# 4) Insert a new bullet after the "Built a loyal customer base..."
#    line.
# ---------------------------------------------------------------------
$loyalPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "for business thrived\.") {
        $loyalPara = $p
        break
    }
}
$loyalPara.Range.InsertParagraphAfter()
$navPara = $loyalPara.Next()
$navPara.Range.InsertBefore("Navigated business landscapes with decisive action, turning any obstacle into opportunities.")

Write-Output "done"
